$d = $word.ActiveDocument

# --- Step 0: capture paragraph 3 ("...年6月2日星期四") and paragraph 4's
# original text ("...初四，明天就是端午节了。") before making any edits, so the
# duplication below is unambiguous.
$p3 = $d.Paragraphs.Item(3)
$srcRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$srcFormattedText = $srcRange.FormattedText

# --- Step 1: update the original (4th/last) paragraph's weather text in place
# first, while it is still the only paragraph containing this text.
$r = $d.Content
$r.Find.Execute("中雨，今天是农历五月初四，明天就是端午节了。", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "中雨，今天是农历五月初五，中国传统端午节。", 2)

# --- Step 2: insert a new paragraph after paragraph 3 that holds the text the
# old paragraph 4 used to have ("...初四，明天就是端午节了。").
$p3.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item(4)
$ip1 = $d.Range($newPara1.Range.End - 1, $newPara1.Range.End - 1)
$ip1.Text = "中雨，今天是农历五月初四，明天就是端午节了。"

# --- Step 3: insert another new paragraph after that one, for "2022年6月3日星期五",
# replicating the exact 3-run split ("2" / "022" / "年6月3日星期五") used by the
# existing date paragraphs (copy paragraph 3's run/format structure, then edit text).
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(5)
$ip2 = $d.Range($newPara2.Range.End - 1, $newPara2.Range.End - 1)
$ip2.Text = "PLACEHOLDER"
$target = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
$target.FormattedText = $srcFormattedText

$pStart = $newPara2.Range.Start
$pEnd = $newPara2.Range.End
$suffixRange = $d.Range($pStart + 4, $pEnd - 1)
$suffixRange.Text = "年6月3日星期五"
